$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 17
$ws.Cells.Item(2, 4).Value = 44259
$ws.Cells.Item(2, 8).Value = 'Sin especificar'
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 4000
$ws.Cells.Item(2, 12).Value = 4500
$ws.Cells.Item(2, 13).Value = 4250
$ws.Cells.Item(2, 16).Value = 4250

# Row 3 <- original row 13
$ws.Cells.Item(3, 4).Value = 44575
$ws.Cells.Item(3, 8).Value = 'Sin especificar'
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 160
$ws.Cells.Item(3, 11).Value = 6500
$ws.Cells.Item(3, 12).Value = 7000
$ws.Cells.Item(3, 13).Value = 6750
$ws.Cells.Item(3, 16).Value = 6750

# Row 4 <- original row 16
$ws.Cells.Item(4, 4).Value = 44789
$ws.Cells.Item(4, 8).Value = 'Sin especificar'
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 5000
$ws.Cells.Item(4, 12).Value = 6000
$ws.Cells.Item(4, 13).Value = 5500
$ws.Cells.Item(4, 16).Value = 5500

# Row 5 <- original row 3
$ws.Cells.Item(5, 4).Value = 44539
$ws.Cells.Item(5, 8).Value = 'Americana (o)'
$ws.Cells.Item(5, 9).Value = 'Primera'
$ws.Cells.Item(5, 10).Value = 160
$ws.Cells.Item(5, 11).Value = 6500
$ws.Cells.Item(5, 12).Value = 7000
$ws.Cells.Item(5, 13).Value = 6750
$ws.Cells.Item(5, 16).Value = 6750

# Row 6 <- original row 15
$ws.Cells.Item(6, 4).Value = 44497
$ws.Cells.Item(6, 8).Value = 'Sin especificar'
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(6, 11).Value = 5000
$ws.Cells.Item(6, 12).Value = 6000
$ws.Cells.Item(6, 13).Value = 5500
$ws.Cells.Item(6, 16).Value = 5500

# Row 7 <- original row 6
$ws.Cells.Item(7, 4).Value = 44253
$ws.Cells.Item(7, 8).Value = 'Americana (o)'
$ws.Cells.Item(7, 9).Value = 'Segunda'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 4000
$ws.Cells.Item(7, 12).Value = 4500
$ws.Cells.Item(7, 13).Value = 4250
$ws.Cells.Item(7, 16).Value = 4250

# Row 8 <- original row 2
$ws.Cells.Item(8, 4).Value = 44699
$ws.Cells.Item(8, 8).Value = 'Sin especificar'
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 50
$ws.Cells.Item(8, 11).Value = 9000
$ws.Cells.Item(8, 12).Value = 9500
$ws.Cells.Item(8, 13).Value = 9250
$ws.Cells.Item(8, 16).Value = 9250

# Row 9 <- original row 8
$ws.Cells.Item(9, 4).Value = 44804
$ws.Cells.Item(9, 8).Value = 'Sin especificar'
$ws.Cells.Item(9, 9).Value = 'Primera'
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 5500
$ws.Cells.Item(9, 12).Value = 6000
$ws.Cells.Item(9, 13).Value = 5750
$ws.Cells.Item(9, 16).Value = 5750

# Row 10 <- original row 4
$ws.Cells.Item(10, 4).Value = 44371
$ws.Cells.Item(10, 8).Value = 'Sin especificar'
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 8000
$ws.Cells.Item(10, 13).Value = 7375
$ws.Cells.Item(10, 16).Value = 7375

# Row 11 <- original row 14
$ws.Cells.Item(11, 4).Value = 44764
$ws.Cells.Item(11, 8).Value = 'Americana (o)'
$ws.Cells.Item(11, 9).Value = 'Primera'
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 7000
$ws.Cells.Item(11, 12).Value = 8000
$ws.Cells.Item(11, 13).Value = 7500
$ws.Cells.Item(11, 16).Value = 7500

# Row 12 <- original row 10
$ws.Cells.Item(12, 4).Value = 44309
$ws.Cells.Item(12, 8).Value = 'Sin especificar'
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 50
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 9000
$ws.Cells.Item(12, 13).Value = 8500
$ws.Cells.Item(12, 16).Value = 8500

# Row 13 <- original row 11
$ws.Cells.Item(13, 4).Value = 44410
$ws.Cells.Item(13, 8).Value = 'Sin especificar'
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 5500
$ws.Cells.Item(13, 12).Value = 6000
$ws.Cells.Item(13, 13).Value = 5750
$ws.Cells.Item(13, 16).Value = 5750

# Row 14 <- original row 7
$ws.Cells.Item(14, 4).Value = 44636
$ws.Cells.Item(14, 8).Value = 'Americana (o)'
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 60
$ws.Cells.Item(14, 11).Value = 8000
$ws.Cells.Item(14, 12).Value = 9000
$ws.Cells.Item(14, 13).Value = 8500
$ws.Cells.Item(14, 16).Value = 8500

# Row 15 <- original row 9
$ws.Cells.Item(15, 4).Value = 44263
$ws.Cells.Item(15, 8).Value = 'Sin especificar'
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 8000
$ws.Cells.Item(15, 13).Value = 7500
$ws.Cells.Item(15, 16).Value = 7500

# Row 16 <- original row 18
$ws.Cells.Item(16, 4).Value = 44281
$ws.Cells.Item(16, 8).Value = 'Sin especificar'
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 100
$ws.Cells.Item(16, 11).Value = 5000
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 13).Value = 5500
$ws.Cells.Item(16, 16).Value = 5500

# Row 17 <- original row 12
$ws.Cells.Item(17, 4).Value = 44559
$ws.Cells.Item(17, 8).Value = 'Americana (o)'
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 100
$ws.Cells.Item(17, 11).Value = 5000
$ws.Cells.Item(17, 12).Value = 6000
$ws.Cells.Item(17, 13).Value = 5500
$ws.Cells.Item(17, 16).Value = 5500

# Row 18 <- original row 5
$ws.Cells.Item(18, 4).Value = 44414
$ws.Cells.Item(18, 8).Value = 'Sin especificar'
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 6000
$ws.Cells.Item(18, 12).Value = 7000
$ws.Cells.Item(18, 13).Value = 6500
$ws.Cells.Item(18, 16).Value = 6500
